$wb = $excel.ActiveWorkbook

# Add the new worksheet right after the last existing sheet ("ODI Bowling")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Match the page margins used by the other sheets in this workbook
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row (bold, bordered, centered like the other sheets' header rows)
$headerRange = $ws.Range("A1:F1")
$ws.Cells.Item(1,1).Value = "MATCH_CODE"
$ws.Cells.Item(1,2).Value = "BATTING_POSITION"
$ws.Cells.Item(1,3).Value = "NUM_4"
$ws.Cells.Item(1,4).Value = "NUM_6"
$ws.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1,6).Value = "MAN_OF_MATCH"
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Row 2 -- MATCH_CODE 4524, rest blank except MAN_OF_MATCH = NO
$ws.Cells.Item(2,1).Value = "'4524"
$c = $ws.Cells.Item(2,2); $c.NumberFormat = "@"; $c.Value = ""
$c = $ws.Cells.Item(2,3); $c.NumberFormat = "@"; $c.Value = ""
$c = $ws.Cells.Item(2,4); $c.NumberFormat = "@"; $c.Value = ""
$c = $ws.Cells.Item(2,5); $c.NumberFormat = "@"; $c.Value = ""
$ws.Cells.Item(2,6).Value = "NO"

# Row 3 -- MATCH_CODE 4526, with numeric batting position and text-typed stats
$ws.Cells.Item(3,1).Value = "'4526"
$ws.Cells.Item(3,2).Value = 6
$ws.Cells.Item(3,3).Value = "'0"
$ws.Cells.Item(3,4).Value = "'1"
$ws.Cells.Item(3,5).Value = "'7.67%"
$ws.Cells.Item(3,6).Value = "NO"
